# Move FAOSTAT extraction scripts into the scenario_analysis_v2 pipeline.
# Two new rows are inserted right before the existing "foreign analysis"
# section (old row 16) to document extract_faostat.R / extract_fbs.R;
# every row that used to start at row 16 shifts down by two (their
# "Order to run" numbers shift the same way, automatically, because they
# come along for the ride with the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the old row 16 ("fao_foreign_land_imports_DT.R"),
# pushing the whole foreign-analysis + visualizations block down to rows 18-25.
$ws.Rows("16:17").Insert()

# Fill the two new rows column-by-column (matches how the sheet was authored).
$ws.Range("A16").Value = "foreign analysis"
$ws.Range("A17").Value = "foreign analysis"

$ws.Range("B16").Value = 15
$ws.Range("B17").Value = 16

$ws.Range("C16").Value = "extract_faostat.R"
$ws.Range("C17").Value = "extract_fbs.R"

$ws.Range("D16").Value = "process raw FAOSTAT data into usable format, averaging the past 5 years of data for all data points"
$ws.Range("D17").Value = "process the food balance sheet data from FAOSTAT, averaging the past 5 years of data for all data points"

$ws.Range("E16").Value = "files in raw_data/FAOSTAT/31aug2020"
$ws.Range("E17").Value = "raw_data/FAOSTAT/31aug2020/FoodBalanceSheets_E_All_Data_(Normalized).csv; raw_data/FAOSTAT/faostat_item_group_lookup.csv"

$ws.Range("F16").Value = "raw FAOSTAT production, emissions, land use, trade, and production value CSVs"
$ws.Range("F17").Value = "raw food balance sheet CSV; lookup table for FAOSTAT item groups showing which codes are aggregated from other codes"

$ws.Range("G16").Value = "faostat2017/(many files).csv"
$ws.Range("G17").Value = "fao_fbs/fbs_*.csv"

$ws.Range("H16").Value = "Processed FAOSTAT data files"
$ws.Range("H17").Value = "Processed FAOSTAT food balance sheet data files"

# Row heights for the two new rows (wrapped-text autosize equivalents).
$ws.Rows("16").RowHeight = 45
$ws.Rows("17").RowHeight = 75

# Leave the active selection on the newly inserted area, like the author did.
$ws.Range("H18").Select()
